$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update existing row 2 values (stats refresh after latest release)
$ws.Range("C2").Value = 224
$ws.Range("D2").Value = 177
$ws.Range("E2").Value = 4368
$ws.Range("F2").Value = 1523
$ws.Range("G2").Value = 3195
$ws.Range("H2").Value = 57442
$ws.Range("I2").Value = 40634
$ws.Range("J2").Value = 1
$ws.Range("L2").Value = 203
$ws.Range("P2").Value = 122
$ws.Range("R2").Value = 541

# Add a new empty row 3 (part of the table range, same styles as row2)
$ws.Range("A2:R2").Copy()
$ws.Range("A3:R3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A3:R3").Value = $null
$excel.CutCopyMode = 0

# Expand the table (ListObject) to include the new row
$table = $ws.ListObjects.Item("Data")
$table.Resize($ws.Range("A1:R3"))

# Update the active selection to A3
$ws.Range("A3").Select()
